# Dragon Spin doc restructure:
#  1. Remove the "Meta description: ..." paragraph that currently sits
#     right after the H1 title.
#  2. At the very end of the document (just before the last "Prompt: ..."
#     paragraph), insert a new bold paragraph carrying the page title
#     text ("Play Dragon Spin Free and Experience the Impressive Sound
#     Design").
#  3. Change the text of the final (previously "Prompt: ...") italic
#     paragraph to the meta-description copy ("Discover Dragon Spin, a
#     low-volatility slot game ...").

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (paragraph #2) ---
$metaParaText = "Meta description: Discover Dragon Spin, a low-volatility slot game with free spins and various bonuses. Try it for free and enjoy the top-notch sound and graphics design."
$metaPara = $d.Paragraphs(2)
if ($metaPara.Range.Text.TrimEnd([char]13) -eq $metaParaText) {
    $metaPara.Range.Delete()
}

# --- Step 2: insert a new bold paragraph right before the last paragraph ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dragon Spin Free and Experience the Impressive Sound Design</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $insertPoint.InsertXML($xmlSnippet)

# InsertXML leaves a stray empty paragraph behind the inserted one (it
# supplies the paragraph break that separates the new paragraph from the
# following "Prompt: ..." paragraph) -- remove that extra paragraph mark
# so the "Prompt: ..." paragraph directly follows the new bold one again.
$newBoldPara = $lastPara.Previous
$strayPara = $newBoldPara.Next
$strayPara.Range.Delete()

# --- Step 3: swap the text of the final paragraph (was "Prompt: ...") ---
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$oldPromptText = 'Prompt: Design a cartoon-style image for the game "Dragon Spin" featuring a happy Maya warrior with glasses. The image should be eye-catching and vibrant, showcasing the mythical dragon theme of the game while also highlighting the fun and playful nature of the Maya warrior character. The warrior should be depicted with a large smile on their face, holding a staff or a sword and standing in a powerful pose. The background of the image should feature a cityscape inspired by ancient Mayan architecture, with a dragon flying in the distance. Color scheme should be vibrant and bold, incorporating shades of red, blue, and yellow. Overall, the image should capture the spirit of adventure and excitement that players can expect when playing "Dragon Spin."'
$newDescText = "Discover Dragon Spin, a low-volatility slot game with free spins and various bonuses. Try it for free and enjoy the top-notch sound and graphics design."

$finalPara.Range.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, $true, 1, $false, $newDescText, 2)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
